$d = $word.ActiveDocument

$replacements = @(
    @("637×2=", "523×5="),
    @("947×3=", "316×8="),
    @("266×3=", "534×9="),
    @("203×3=", "800×8="),
    @("710×6=", "554×7="),
    @("375×7=", "518×2="),
    @("159×7=", "875×5="),
    @("145×5=", "571×5="),
    @("928×9=", "951×6="),
    @("891×8=", "478×5="),
    @("141×5=", "685×5="),
    @("510×4=", "442×7="),
    @("546×5=", "887×4="),
    @("683×9=", "387×2="),
    @("521×2=", "206×5="),
    @("566×4=", "522×4="),
    @("186×4=", "442×4="),
    @("439×7=", "525×6="),
    @("254×5=", "726×7="),
    @("989×3=", "601×3="),
    @("848×9=", "219×5="),
    @("649×9=", "297×7="),
    @("705×9=", "414×8="),
    @("144×9=", "102×3="),
    @("918×7=", "898×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
